# Commit: "Add data for 2022-03-04"
# Updates the "carjacking by neighborhood by month" report to the next day's
# cut: renames the sheet/header from "through February 23" to "through
# February 24", inserts the newly-appearing "Armour Square" neighborhood row
# (alphabetically between Andersonville and Avondale), and bumps a handful of
# month-column counts (including the brand new row) to reflect the extra
# day's incidents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-02-24"

# 2) Update the matching header label in A1's row (current-month column).
$ws.Range("B1").Value = "February 2022 (through February 24)"

# 3) Insert the new "Armour Square" row. It belongs alphabetically right
#    after "Andersonville" (row 58) and before "Avondale" (currently row 59),
#    so insert a blank row at 59 and push Avondale..Wrigleyville down by one.
$ws.Rows(59).Insert()

# Copy the formatting from the row above (Andersonville) so the new label
# cell picks up the same bold/centered/bordered style used by every other
# neighborhood-name cell in column A.
$ws.Range("A58").Copy()
$ws.Range("A59").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A59").Value = "Armour Square"
$ws.Range("F59").Value = 1

# 4) New incident recorded for Dunning (now row 64 after the insert above).
$ws.Range("N64").Value = 1

# 5) Updated counts for neighborhoods unaffected by the row insertion
#    (all above row 58).
$ws.Range("D3").Value = 13    # Austin
$ws.Range("L5").Value = 1     # Calumet Heights
$ws.Range("B8").Value = 7     # North Lawndale
$ws.Range("B11").Value = 8    # Garfield Park
$ws.Range("L11").Value = 2    # Garfield Park
$ws.Range("N11").Value = 2    # Garfield Park
$ws.Range("F14").Value = 1    # Bridgeport
$ws.Range("D15").Value = 5    # West Town
$ws.Range("L16").Value = 1    # Chicago Lawn
$ws.Range("B20").Value = 4    # Kenwood
$ws.Range("D22").Value = 4    # Humboldt Park
$ws.Range("L23").Value = 1    # Washington Heights
$ws.Range("D49").Value = 1    # Irving Park
$ws.Range("N54").Value = 1    # Belmont Cragin
